# Expand the "Publikus IP-címek" worksheet with the full, current IP address
# inventory: updated router-interface labels, a new "Távmunkás SOHO" row,
# two new sub-tables (WAN/LAN for the remote-worker setup, and the private
# IP ranges used at each site / secretary office), plus a light accent fill
# on the header rows and section dividers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell values
# ---------------------------------------------------------------------

# Row 1 - top banner (unchanged text)
$ws.Range("A1").Value = 'Publikus IP-címek az ISP felé'
$ws.Range("B1").Value = 'kábelkötések'
$ws.Range("C1").Value = '4-esével növekszik'

# Row 2 - column headers (C2 label clarified)
$ws.Range("A2").Value = 'Telephely'
$ws.Range("B2").Value = 'Router interfész (ISP felé)'
$ws.Range("C2").Value = 'Publikus IP (Routeré, az isp routeren 1-el kisebb szám)'

# Rows 3-5 - existing sites, interface names updated to match the current topology
$ws.Range("A3").Value = 'Budapest - R1'
$ws.Range("B3").Value = 'G0/2 → ISP_R G0/1'
$ws.Range("C3").Value = '203.0.113.2/30'

$ws.Range("A4").Value = 'Szeged - R2'
$ws.Range("B4").Value = 'G0/2 → ISP_R G0/2'
$ws.Range("C4").Value = '203.0.113.6/30'

$ws.Range("A5").Value = 'Debrecen - R3'
$ws.Range("B5").Value = 'Gig0/0/0(optikai kábel itt) → Gig0/2/0'
$ws.Range("C5").Value = '203.0.113.10/30'

# Row 6 - new site: the remote-worker SOHO router
$ws.Range("A6").Value = 'Távmunkás SOHO'
$ws.Range("B6").Value = 'ROUTER_INTERFACE → ISP_R G0/0'
$ws.Range("C6").Value = '1.1.1.1/24'

# Row 7 - section divider, renamed to mention NAT
$ws.Range("A7").Value = 'Távmunkás környezet - NAT'
$ws.Range("B7").Value = ''
$ws.Range("C7").Value = ''

# Rows 8-9 - WAN/LAN description for the remote-worker environment
$ws.Range("A8").Value = 'WAN - 1.1.1.0 255.255.255.0/24'
$ws.Range("A9").Value = 'LAN - 192.168.50.0 255.255.255.0/24 vlan 20-ban működik'

# Row 11 - new section: private IP ranges per site / secretary office
$ws.Range("A11").Value = 'Privát IP-k a telephelyeken (titkárságok egy másik közeli épületben)'
$ws.Range("B11").Value = ''
$ws.Range("C11").Value = ''

$ws.Range("A12").Value = 'Telephely'
$ws.Range("B12").Value = 'Hálózat'
$ws.Range("C12").Value = 'Router IP'

$ws.Range("A13").Value = 'Budapest (központ)'
$ws.Range("B13").Value = '192.168.10.0/24'
$ws.Range("C13").Value = '192.168.10.1'

$ws.Range("A14").Value = 'Debrecen'
$ws.Range("B14").Value = '192.168.20.0/24'
$ws.Range("C14").Value = '192.168.20.1'

$ws.Range("A15").Value = 'Szeged'
$ws.Range("B15").Value = '192.168.30.0/24'
$ws.Range("C15").Value = '192.168.30.1'

$ws.Range("A16").Value = 'Budapest titkárság'
$ws.Range("B16").Value = '192.168.11.0/24'
$ws.Range("C16").Value = '192.168.11.1'

$ws.Range("A17").Value = 'Szeged titkárság'
$ws.Range("B17").Value = '192.168.31.0/24'
$ws.Range("C17").Value = '192.168.31.1'

# ---------------------------------------------------------------------
# 2. Column widths / row heights
# ---------------------------------------------------------------------

# Column B needs to be a bit wider for the longer interface descriptions
$ws.Columns.Item(2).ColumnWidth = 13.6

$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 45

# ---------------------------------------------------------------------
# 3. Formatting
# ---------------------------------------------------------------------

$accentTheme = 10
$accentTint = 0.59996337778862885

# Banner row (1) and the two new blank divider cells - accent fill only
$fillOnly = $ws.Range("A1:C1,B7:C7,B11:C11")
$fillOnly.Interior.ThemeColor = $accentTheme
$fillOnly.Interior.TintAndShade = $accentTint

# Column header row (2) - bold, centered, accent fill
$headerRow = $ws.Range("A2:C2")
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4108
$headerRow.WrapText = $true
$headerRow.Interior.ThemeColor = $accentTheme
$headerRow.Interior.TintAndShade = $accentTint

# Section-divider label cells (7, 11) - accent fill, wrapped, vertically centered
$dividers = $ws.Range("A7,A11")
$dividers.VerticalAlignment = -4108
$dividers.WrapText = $true
$dividers.Interior.ThemeColor = $accentTheme
$dividers.Interior.TintAndShade = $accentTint

# Plain wrapped body cells (site/interface names)
$bodyA = $ws.Range("A3:B6,A8:A9")
$bodyA.VerticalAlignment = -4108
$bodyA.WrapText = $true

# IP-address column uses the narrow Arial Unicode MS font
$ipCol = $ws.Range("C3:C6")
$ipCol.Font.Name = 'Arial Unicode MS'
$ipCol.Font.Size = 10
$ipCol.VerticalAlignment = -4108
$ipCol.WrapText = $true

# ---------------------------------------------------------------------
# 4. Selection (matches where the author left the cursor)
# ---------------------------------------------------------------------
$ws.Range("G4").Select()
